# Auto-generated edit script: updates market-price derived columns (H-N)
# across multiple sheets to match refreshed Universalis price data.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 179.53847
$ws.Cells.Item(6, 9).Value = 190.44444
$ws.Cells.Item(6, 11).Value = 571.33332
$ws.Cells.Item(6, 13).Value = -459.33332
$ws.Cells.Item(38, 8).Value = 1316.5714
$ws.Cells.Item(38, 9).Value = 353.2
$ws.Cells.Item(38, 11).Value = 1059.6
$ws.Cells.Item(38, 13).Value = -687.5999999999999
$ws.Cells.Item(40, 8).Value = 2204.2083
$ws.Cells.Item(40, 9).Value = 1833.1666
$ws.Cells.Item(40, 10).Value = 2575.25
$ws.Cells.Item(40, 11).Value = 1833.1666
$ws.Cells.Item(40, 12).Value = 2575.25
$ws.Cells.Item(40, 13).Value = -1658.1666
$ws.Cells.Item(40, 14).Value = -2925.25
$ws.Cells.Item(70, 8).Value = 88520.75
$ws.Cells.Item(70, 9).Value = 1900
$ws.Cells.Item(70, 10).Value = 175141.5
$ws.Cells.Item(70, 11).Value = 5700
$ws.Cells.Item(70, 12).Value = 525424.5
$ws.Cells.Item(70, 13).Value = -5430
$ws.Cells.Item(70, 14).Value = -525964.5
$ws.Cells.Item(73, 8).Value = 88520.75
$ws.Cells.Item(73, 9).Value = 1900
$ws.Cells.Item(73, 10).Value = 175141.5
$ws.Cells.Item(73, 11).Value = 5700
$ws.Cells.Item(73, 12).Value = 525424.5
$ws.Cells.Item(73, 13).Value = -4764
$ws.Cells.Item(73, 14).Value = -527296.5
$ws.Cells.Item(80, 8).Value = 377.33334
$ws.Cells.Item(80, 9).Value = 267.5
$ws.Cells.Item(80, 10).Value = 432.25
$ws.Cells.Item(80, 11).Value = 802.5
$ws.Cells.Item(80, 12).Value = 1296.75
$ws.Cells.Item(80, 13).Value = 195.5
$ws.Cells.Item(80, 14).Value = -3292.75
$ws.Cells.Item(83, 8).Value = 377.33334
$ws.Cells.Item(83, 9).Value = 267.5
$ws.Cells.Item(83, 10).Value = 432.25
$ws.Cells.Item(83, 11).Value = 2407.5
$ws.Cells.Item(83, 12).Value = 3890.25
$ws.Cells.Item(83, 13).Value = 2584.5
$ws.Cells.Item(83, 14).Value = -13874.25
$ws.Cells.Item(88, 8).Value = 2124.75
$ws.Cells.Item(88, 10).Value = 3500
$ws.Cells.Item(88, 12).Value = 3500
$ws.Cells.Item(88, 14).Value = -4312
$ws.Cells.Item(91, 8).Value = 2124.75
$ws.Cells.Item(91, 10).Value = 3500
$ws.Cells.Item(91, 12).Value = 3500
$ws.Cells.Item(91, 14).Value = -6308
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 9).Value = 0
$ws.Cells.Item(92, 11).Value = 0
$ws.Cells.Item(92, 13).ClearContents()
$ws.Cells.Item(98, 8).Value = 3064.2222
$ws.Cells.Item(98, 9).Value = 2321.5
$ws.Cells.Item(98, 11).Value = 2321.5
$ws.Cells.Item(98, 13).Value = -823.5
$ws.Cells.Item(103, 8).Value = 4475
$ws.Cells.Item(103, 9).Value = 2900
$ws.Cells.Item(103, 11).Value = 8700
$ws.Cells.Item(103, 13).Value = -8114
$ws.Cells.Item(122, 8).Value = 3064.2222
$ws.Cells.Item(122, 9).Value = 2321.5
$ws.Cells.Item(122, 11).Value = 6964.5
$ws.Cells.Item(122, 13).Value = -4514.5
$ws.Cells.Item(132, 8).Value = 2609.6
$ws.Cells.Item(132, 9).Value = 2609.6
$ws.Cells.Item(132, 11).Value = 7828.799999999999
$ws.Cells.Item(132, 13).Value = -5298.799999999999
$ws.Cells.Item(137, 8).Value = 1588.3636
$ws.Cells.Item(137, 9).Value = 1434.3529
$ws.Cells.Item(137, 10).Value = 2112
$ws.Cells.Item(137, 11).Value = 4303.0587
$ws.Cells.Item(137, 12).Value = 6336
$ws.Cells.Item(137, 13).Value = -1753.0587
$ws.Cells.Item(137, 14).Value = -11436
$ws.Cells.Item(138, 8).Value = 1815.3334
$ws.Cells.Item(138, 9).Value = 1815.3334
$ws.Cells.Item(138, 11).Value = 5446.0002
$ws.Cells.Item(138, 13).Value = -306.0002000000004

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1236.75
$ws.Cells.Item(61, 9).Value = 982.3333
$ws.Cells.Item(61, 11).Value = 982.3333
$ws.Cells.Item(61, 13).Value = -770.3333
$ws.Cells.Item(132, 8).Value = 1395.6666
$ws.Cells.Item(132, 9).Value = 1320.125
$ws.Cells.Item(132, 11).Value = 3960.375
$ws.Cells.Item(132, 13).Value = -1430.375
$ws.Cells.Item(136, 8).Value = 1236.75
$ws.Cells.Item(136, 9).Value = 982.3333
$ws.Cells.Item(136, 11).Value = 2946.9999
$ws.Cells.Item(136, 13).Value = -396.9998999999998

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1058.3549
$ws.Cells.Item(94, 9).Value = 552
$ws.Cells.Item(94, 11).Value = 552
$ws.Cells.Item(94, 13).Value = -101

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 68817.336
$ws.Cells.Item(62, 10).Value = 134833
$ws.Cells.Item(62, 12).Value = 134833
$ws.Cells.Item(62, 14).Value = -136081
$ws.Cells.Item(65, 8).Value = 68817.336
$ws.Cells.Item(65, 10).Value = 134833
$ws.Cells.Item(65, 12).Value = 674165
$ws.Cells.Item(65, 14).Value = -680405
$ws.Cells.Item(134, 8).Value = 2851.2307
$ws.Cells.Item(134, 9).Value = 2454.8572
$ws.Cells.Item(134, 10).Value = 3313.6667
$ws.Cells.Item(134, 11).Value = 7364.571599999999
$ws.Cells.Item(134, 12).Value = 9941.000100000001
$ws.Cells.Item(134, 13).Value = -4829.571599999999
$ws.Cells.Item(134, 14).Value = -15011.0001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 607.7059
$ws.Cells.Item(4, 9).Value = 559.4138
$ws.Cells.Item(4, 11).Value = 1678.2414
$ws.Cells.Item(4, 13).Value = -1566.2414
$ws.Cells.Item(17, 8).Value = 1185
$ws.Cells.Item(17, 10).Value = 1841.6666
$ws.Cells.Item(17, 12).Value = 5524.9998
$ws.Cells.Item(17, 14).Value = -5862.9998
$ws.Cells.Item(34, 8).Value = 1722.3
$ws.Cells.Item(34, 10).Value = 2759.6
$ws.Cells.Item(34, 12).Value = 8278.799999999999
$ws.Cells.Item(34, 14).Value = -8446.799999999999
$ws.Cells.Item(36, 8).Value = 1612.5
$ws.Cells.Item(36, 9).Value = 1612.5
$ws.Cells.Item(36, 10).Value = 0
$ws.Cells.Item(36, 11).Value = 4837.5
$ws.Cells.Item(36, 12).Value = 0
$ws.Cells.Item(36, 13).ClearContents()
$ws.Cells.Item(36, 14).Value = -4668.5
$ws.Cells.Item(55, 8).Value = 147657
$ws.Cells.Item(55, 10).Value = 6699.8
$ws.Cells.Item(55, 12).Value = 20099.4
$ws.Cells.Item(55, 14).Value = -20453.4
$ws.Cells.Item(121, 8).Value = 449.8889
$ws.Cells.Item(121, 10).Value = 535
$ws.Cells.Item(121, 12).Value = 1605
$ws.Cells.Item(121, 14).Value = -4225
$ws.Cells.Item(137, 8).Value = 4753.3335
$ws.Cells.Item(137, 10).Value = 4350
$ws.Cells.Item(137, 12).Value = 13050
$ws.Cells.Item(137, 14).Value = -23250

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 4335.857
$ws.Cells.Item(80, 9).Value = 3781.125
$ws.Cells.Item(80, 10).Value = 4677.231
$ws.Cells.Item(80, 11).Value = 3781.125
$ws.Cells.Item(80, 12).Value = 4677.231
$ws.Cells.Item(80, 13).Value = -2783.125
$ws.Cells.Item(80, 14).Value = -6673.231
$ws.Cells.Item(83, 8).Value = 4335.857
$ws.Cells.Item(83, 9).Value = 3781.125
$ws.Cells.Item(83, 10).Value = 4677.231
$ws.Cells.Item(83, 11).Value = 18905.625
$ws.Cells.Item(83, 12).Value = 23386.155
$ws.Cells.Item(83, 13).Value = -13913.625
$ws.Cells.Item(83, 14).Value = -33370.155
$ws.Cells.Item(122, 8).Value = 69224.734
$ws.Cells.Item(122, 9).Value = 2125.182
$ws.Cells.Item(122, 11).Value = 6375.545999999999
$ws.Cells.Item(122, 13).Value = -3925.545999999999
$ws.Cells.Item(123, 8).Value = 42074.918
$ws.Cells.Item(123, 10).Value = 42727.184
$ws.Cells.Item(123, 12).Value = 42727.184
$ws.Cells.Item(123, 14).Value = -47627.184
$ws.Cells.Item(132, 8).Value = 1789.7368
$ws.Cells.Item(132, 9).Value = 1486.7858
$ws.Cells.Item(132, 11).Value = 4460.357400000001
$ws.Cells.Item(132, 13).Value = -1930.357400000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 700
$ws.Cells.Item(7, 9).Value = 700
$ws.Cells.Item(7, 11).Value = 700
$ws.Cells.Item(7, 13).Value = -588
$ws.Cells.Item(40, 8).Value = 2584.5
$ws.Cells.Item(40, 9).Value = 2584.5
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 2584.5
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 13).ClearContents()
$ws.Cells.Item(40, 14).Value = -2448.5
$ws.Cells.Item(55, 8).Value = 1495.75
$ws.Cells.Item(55, 9).Value = 994
$ws.Cells.Item(55, 10).Value = 1997.5
$ws.Cells.Item(55, 11).Value = 994
$ws.Cells.Item(55, 12).Value = 1997.5
$ws.Cells.Item(55, 13).Value = -821
$ws.Cells.Item(55, 14).Value = -2343.5
$ws.Cells.Item(126, 8).Value = 700
$ws.Cells.Item(126, 9).Value = 700
$ws.Cells.Item(126, 11).Value = 2100
$ws.Cells.Item(126, 13).Value = 370
$ws.Cells.Item(132, 8).Value = 3000
$ws.Cells.Item(132, 9).Value = 3000
$ws.Cells.Item(132, 10).Value = 3000
$ws.Cells.Item(132, 11).Value = 9000
$ws.Cells.Item(132, 12).Value = 9000
$ws.Cells.Item(132, 13).Value = -6470
$ws.Cells.Item(132, 14).Value = -14060
$ws.Cells.Item(133, 8).Value = 59999.5
$ws.Cells.Item(133, 10).Value = 59999.5
$ws.Cells.Item(133, 12).Value = 59999.5
$ws.Cells.Item(133, 14).Value = -65059.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 500
$ws.Cells.Item(2, 10).Value = 500
$ws.Cells.Item(2, 12).Value = 500
$ws.Cells.Item(2, 14).Value = -724
$ws.Cells.Item(42, 8).Value = 19997.5
$ws.Cells.Item(42, 9).Value = 19997.5
$ws.Cells.Item(42, 11).Value = 19997.5
$ws.Cells.Item(42, 13).Value = -19619.5
$ws.Cells.Item(43, 8).Value = 10500.25
$ws.Cells.Item(43, 9).Value = 10500.25
$ws.Cells.Item(43, 11).Value = 10500.25
$ws.Cells.Item(43, 13).Value = -10351.25
$ws.Cells.Item(95, 8).Value = 14172
$ws.Cells.Item(95, 10).Value = 14172
$ws.Cells.Item(95, 12).Value = 14172
$ws.Cells.Item(95, 14).Value = -19664
$ws.Cells.Item(122, 8).Value = 1825.5
$ws.Cells.Item(122, 10).Value = 1799
$ws.Cells.Item(122, 12).Value = 5397
$ws.Cells.Item(122, 14).Value = -10297
$ws.Cells.Item(126, 8).Value = 3000
$ws.Cells.Item(126, 9).Value = 3000
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 9000
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 13).ClearContents()
$ws.Cells.Item(126, 14).Value = -6530
$ws.Cells.Item(136, 8).Value = 2541.25
$ws.Cells.Item(136, 9).Value = 1083.25
$ws.Cells.Item(136, 11).Value = 3249.75
$ws.Cells.Item(136, 13).Value = -699.75
